$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 112 (shifts old rows 112..222 down to 113..223)
$ws.Rows(112).Insert()

# Populate the newly inserted row 112 with the new record
$ws.Cells.Item(112, 1).Value = 3
$ws.Cells.Item(112, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(112, 3).Value = "Coquimbo"
$ws.Cells.Item(112, 4).Value = 44586
$ws.Cells.Item(112, 5).Value = 5
$ws.Cells.Item(112, 6).Value = 100112001
$ws.Cells.Item(112, 7).Value = "Berenjena"
$ws.Cells.Item(112, 8).Value = "Sin especificar"
$ws.Cells.Item(112, 9).Value = "Primera"
$ws.Cells.Item(112, 10).Value = 73
$ws.Cells.Item(112, 11).Value = 9500
$ws.Cells.Item(112, 12).Value = 10000
$ws.Cells.Item(112, 13).Value = 9760
$ws.Cells.Item(112, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(112, 15).Value = "Región Metropolitana"
$ws.Cells.Item(112, 16).Value = 163
$ws.Cells.Item(112, 17).Value = 60
$ws.Cells.Item(112, 18).Value = "Hortaliza"
